$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the existing "sum" header cell (G1) onto the new
# "Save" header cell (H1) so it reuses the same cell style (bold, bordered,
# centered) rather than minting a brand-new style entry.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Header + values for the new "Save" column
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
